# Regenerate the s_vals data to filter save games.
# Updates columns B, C, D, E, G for rows 2-11 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    # row, B, C, D, E, G
    ,@(2,  1.505614041169197,  1.65323645889881,   3.082599426703578,  0.4998867070740569, 6.741336633845642)
    ,@(3,  0.06328177979961902,0.05231270169004087,0.7127328510149897, 0.4998867070740569, 1.328214039578707)
    ,@(4,  3.182878228561681,  1.65323645889881,   16.98373111632243,  0.4998867070740569, 22.31973251085698)
    ,@(5,  3.182878228561681,  1.65323645889881,   0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    ,@(6,  3.182878228561681,  1.65323645889881,   0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    ,@(7,  1.505614041169197,  0.05231270169004087,0.7127328510149897, 0.4998867070740569, 2.770546300948285)
    ,@(8,  3.182878228561681,  1.65323645889881,   0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    ,@(9,  0.7287194209349384, 1.65323645889881,   3.082599426703578,  6.48142807727062,   11.94598338380795)
    ,@(10, 0.02258322285507441,0.3375848360084654, 0.1529057820181812, 6.48142807727062,   6.994501918152341)
    ,@(11, 0.1554434735375247, 0.3375848360084654, 3.082599426703578,  0.4998867070740569, 4.075514443323626)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]  # B
    $ws.Cells.Item($r, 3).Value = $entry[2]  # C
    $ws.Cells.Item($r, 4).Value = $entry[3]  # D
    $ws.Cells.Item($r, 5).Value = $entry[4]  # E
    $ws.Cells.Item($r, 7).Value = $entry[5]  # G
}
